# Government Revenue Accounting.xlsx - "updated rmi files 3.4.3"
#
# 1) Remove the stray date stamp in About!C1 (and its date-format style).
# 2) On "Set Values Here", update the fuel tax revenue row (row 9) weights:
#       C9: 0 -> 5
#       D9: 5 -> 0
#       F9: 0 -> 5
#    (B9 and E9 stay 0). The dependent percentage formulas (row 23) and the
#    per-policy TRANSPOSE sheets (e.g. GRA-fueltax) recalculate automatically.

$wb = $excel.ActiveWorkbook

# --- About sheet: clear the leftover date value in C1 ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Clear()

# --- Set Values Here: fix fuel tax revenue (row 9) weighting ---
$setValues = $wb.Worksheets.Item("Set Values Here")
$setValues.Range("C9").Value = 5
$setValues.Range("D9").Value = 0
$setValues.Range("F9").Value = 5
[void]$setValues.Range("F10").Select()

# Restore "About" as the active sheet/tab (it was selected before editing).
[void]$about.Activate()

$wb.Application.Calculate() | Out-Null
